$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 2, pushing existing data (rows 2-11) down to rows 3-12
$ws.Rows.Item(2).Insert()

# Reset the selection to a single cell, as seen after the edit
$ws.Range("F8").Select()
